# Appends 5 new survey response rows (421-425) to the bottom of the sheet,
# matching the same shape/format as the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 421 ---------------------------------------------------------------
$ws.Range("A421").Value = 44246.68990465278
$ws.Range("A421").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B421").Value = "PK-Seutu"
$ws.Range("C421").Value = "31-35 v"
$ws.Range("D421").Value = "mies"
$ws.Range("E421").Value = 11
$ws.Range("F421").Value = "Työntekijä / palkollinen"
$ws.Range("G421").Value = 1
$ws.Range("H421").Value = "Full stack"
$ws.Range("I421").Value = "50/50"
$ws.Range("J421").Value = 7000
$ws.Range("K421").Value = 87500
$ws.Range("L421").Value = $true
$ws.Range("M421").Value = "Mavericks"

# --- Row 422 ---------------------------------------------------------------
$ws.Range("A422").Value = 44246.69036510416
$ws.Range("A422").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B422").Value = "PK-Seutu"
$ws.Range("C422").Value = "31-35 v"
$ws.Range("D422").Value = "mies"
$ws.Range("E422").Value = 12
$ws.Range("F422").Value = "Työntekijä / palkollinen"
$ws.Range("G422").Value = 1
$ws.Range("H422").Value = "full-stack"
$ws.Range("I422").Value = "Etä"
$ws.Range("J422").Value = 8000
$ws.Range("K422").Value = 95000
$ws.Range("L422").Value = $true
$ws.Range("M422").Value = "Mavericks"

# --- Row 423 ---------------------------------------------------------------
$ws.Range("A423").Value = 44246.69231409722
$ws.Range("A423").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B423").Value = "Tampere"
$ws.Range("C423").Value = "41-45 v"
$ws.Range("D423").Value = "mies"
$ws.Range("E423").Value = 22
$ws.Range("F423").Value = "Työntekijä / palkollinen"
$ws.Range("G423").Value = 0.8
$ws.Range("H423").Value = "ohjelmistokehittäjä (backend) / arkkitehti"
$ws.Range("I423").Value = "Etä"
$ws.Range("J423").Value = 4700
$ws.Range("K423").Value = 58750
$ws.Range("L423").Value = $false

# --- Row 424 ---------------------------------------------------------------
$ws.Range("A424").Value = 44246.69353475695
$ws.Range("A424").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B424").Value = "PK-Seutu"
$ws.Range("C424").Value = "36-40 v"
$ws.Range("D424").Value = "mies"
$ws.Range("E424").Value = 2
$ws.Range("F424").Value = "Työntekijä / palkollinen"
$ws.Range("G424").Value = 1
$ws.Range("H424").Value = "WordPress-kehittäjä"
$ws.Range("I424").Value = "50/50"
$ws.Range("J424").Value = 3000
$ws.Range("K424").Value = 37500
$ws.Range("L424").Value = $false

# --- Row 425 ---------------------------------------------------------------
$ws.Range("A425").Value = 44246.69392165509
$ws.Range("A425").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B425").Value = "Tampere"
$ws.Range("C425").Value = "31-35 v"
$ws.Range("D425").Value = "mies"
$ws.Range("E425").Value = 5
$ws.Range("F425").Value = "Työntekijä / palkollinen"
$ws.Range("G425").Value = 1
$ws.Range("H425").Value = "Data scientist"
$ws.Range("I425").Value = "Etä"
$ws.Range("J425").Value = 4300
$ws.Range("K425").Value = 53750
$ws.Range("M425").Value = "Wapice"
